# Weekly update: a new price observation for "Arveja Verde" at Terminal
# Hortofrutícola Agro Chillán was recorded. This inserts a new data row at
# row 19 (just below the most-recent existing entry in row 18), pushing the
# previously-existing rows 19-69 down to rows 20-70, and fills the new row
# with the freshly reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 19 - everything below (old rows 19-69)
# shifts down to rows 20-70, preserving all of their data untouched.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new observation.
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44623
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112022
$ws.Range("G19").Value = "Arveja Verde"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 24000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 24500
$ws.Range("N19").Value = "`$/saco 25 kilos"
$ws.Range("O19").Value = "Provincia de Diguillín"
$ws.Range("P19").Value = 980
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
